# Ft: addUsersValidation email is required input
#
# Adds an "email" column of hyperlinked mailto addresses for the two rows
# that previously had no value in column G (row 2 -> James, row 4 -> Jack),
# mirroring the existing "role" entry already present for row 3 (Grant).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add james@gmail.com as a mailto hyperlink in G2
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:james@gmail.com", [Type]::Missing, [Type]::Missing, "james@gmail.com")

# Row 4: add jack@gmail.com as a mailto hyperlink in G4
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:jack@gmail.com", [Type]::Missing, [Type]::Missing, "jack@gmail.com")

# Leave the final selection on G4, matching the saved workbook state
[void]$ws.Range("G4").Select()
